# Add a new data row (row 28) to the "Mixing Block Data" sheet, matching
# the formatting of the preceding row (row 27, which already carries the
# "FAIL" highlight style) and filling in the new reading's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mixing Block Data")

$newRow = 28
$templateRow = 27

# 1) Clone row 27's formatting (fill/border/etc.) onto row 28 so every
#    cell in the new row starts out with the same style as its sibling.
$ws.Range("A" + $templateRow + ":O" + $templateRow).Copy()
$ws.Range("A" + $newRow + ":O" + $newRow).PasteSpecial(-4122)

# 2) Fill in the text/number values for the new row.
$ws.Cells.Item($newRow, 1).Value = "2025-11-21 14:16:06"
$ws.Cells.Item($newRow, 2).Value = "SA01"
$ws.Cells.Item($newRow, 3).Value = "Mixing Block"
$ws.Cells.Item($newRow, 4).Value = "A"
$ws.Cells.Item($newRow, 6).Value = "IN"
$ws.Cells.Item($newRow, 7).Value = "H1"
$ws.Cells.Item($newRow, 8).Value = "Inner"
$ws.Cells.Item($newRow, 9).Value = 111
$ws.Cells.Item($newRow, 10).Value = 4
$ws.Cells.Item($newRow, 11).Value = 3.5
$ws.Cells.Item($newRow, 12).Value = 4.5
$ws.Cells.Item($newRow, 13).Value = "FAIL"

# Column E ("Piece ID") holds the text "4" (not the number 4) in the
# source data, same as rows 26/27. Assigning the string directly would
# get auto-coerced to a number, so instead enter it as a text formula
# and then collapse the formula down to its resulting static text value
# via a values-only paste (this keeps the cell's string type and the
# style applied in step 1 intact).
$ws.Cells.Item($newRow, 5).Formula = '="4"'
$ws.Cells.Item($newRow, 5).Copy()
$ws.Cells.Item($newRow, 5).PasteSpecial(-4163)

$excel.CutCopyMode = 0
